# Update "想去人数" (F column) counts and one cover image URL (I14)
# on both the "展览" and "全部类型" sheets, which hold identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row => new F value (want-to-go count)
$fUpdates = @{
    3  = 1742
    4  = 798
    6  = 1128
    8  = 12082
    12 = 421
    13 = 1117
    14 = 874
    15 = 13514
    16 = 13553
    18 = 156
    19 = 21
    21 = 997
    24 = 1992
    25 = 183
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    $ws.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202403/PlZCFPVs1710502485559.jpeg"
}
